# Apply updated probability values to Sheet1 (team-specific matrix).
# Source: "changes to team matrices from games pulled march 7"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1954314720812183
$ws.Range("C2").Value = 0.550761421319797
$ws.Range("J2").Value = 0.03807106598984772
$ws.Range("P2").Value = 0.149746192893401
$ws.Range("S2").Value = 0.06598984771573604
$ws.Range("B3").Value = 0.01276595744680851
$ws.Range("C3").Value = 0.02127659574468085
$ws.Range("J3").Value = 0.03829787234042553
$ws.Range("P3").Value = 0.7319148936170212
$ws.Range("S3").Value = 0.1957446808510638
$ws.Range("B6").Value = 0.08482142857142858
$ws.Range("E6").Value = 0.004464285714285714
$ws.Range("F6").Value = 0.0625
$ws.Range("J6").Value = 0.3035714285714285
$ws.Range("O6").Value = 0.03571428571428571
$ws.Range("Q6").Value = 0.1651785714285714
$ws.Range("R6").Value = 0.06696428571428571
$ws.Range("S6").Value = 0.2767857142857143
$ws.Range("B7").Value = 0.1058823529411765
$ws.Range("D7").Value = 0.01764705882352941
$ws.Range("F7").Value = 0.04117647058823529
$ws.Range("J7").Value = 0.1941176470588235
$ws.Range("O7").Value = 0.02941176470588235
$ws.Range("Q7").Value = 0.2117647058823529
$ws.Range("R7").Value = 0.05882352941176471
$ws.Range("S7").Value = 0.3411764705882353
$ws.Range("B8").Value = 0.1012658227848101
$ws.Range("D8").Value = 0.01265822784810127
$ws.Range("E8").Value = 0.002531645569620253
$ws.Range("F8").Value = 0.05316455696202532
$ws.Range("J8").Value = 0.1265822784810127
$ws.Range("O8").Value = 0.01012658227848101
$ws.Range("Q8").Value = 0.2126582278481013
$ws.Range("R8").Value = 0.1215189873417721
$ws.Range("S8").Value = 0.3594936708860759
$ws.Range("B9").Value = 0.1557377049180328
$ws.Range("D9").Value = 0.02459016393442623
$ws.Range("F9").Value = 0.04098360655737705
$ws.Range("J9").Value = 0.1311475409836066
$ws.Range("O9").Value = 0.01639344262295082
$ws.Range("Q9").Value = 0.2131147540983606
$ws.Range("R9").Value = 0.1065573770491803
$ws.Range("S9").Value = 0.3114754098360656
$ws.Range("B10").Value = 0.1347424042272127
$ws.Range("D10").Value = 0.02509907529722589
$ws.Range("F10").Value = 0.05812417437252312
$ws.Range("J10").Value = 0.1387054161162483
$ws.Range("O10").Value = 0.02113606340819022
$ws.Range("Q10").Value = 0.2199471598414795
$ws.Range("R10").Value = 0.08586525759577279
$ws.Range("S10").Value = 0.3163804491413474
$ws.Range("G11").Value = 0.1689189189189189
$ws.Range("J11").Value = 0.1047297297297297
$ws.Range("K11").Value = 0.2331081081081081
$ws.Range("L11").Value = 0.4763513513513514
$ws.Range("S11").Value = 0.01689189189189189
$ws.Range("G12").Value = 0.673469387755102
$ws.Range("J12").Value = 0.2585034013605442
$ws.Range("K12").Value = 0.006802721088435374
$ws.Range("L12").Value = 0.02040816326530612
$ws.Range("S12").Value = 0.04081632653061224
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02857142857142857
$ws.Range("H15").Value = 0.1428571428571428
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.4178571428571429
$ws.Range("K15").Value = 0.04642857142857143
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("O15").Value = 0.075
$ws.Range("S15").Value = 0.2214285714285714
$ws.Range("F16").Value = 0.04382470119521913
$ws.Range("H16").Value = 0.1155378486055777
$ws.Range("I16").Value = 0.05577689243027888
$ws.Range("J16").Value = 0.4462151394422311
$ws.Range("K16").Value = 0.1075697211155379
$ws.Range("M16").Value = 0.01195219123505976
$ws.Range("O16").Value = 0.05976095617529881
$ws.Range("S16").Value = 0.1593625498007968
$ws.Range("F17").Value = 0.025
$ws.Range("H17").Value = 0.1269230769230769
$ws.Range("I17").Value = 0.07692307692307693
$ws.Range("J17").Value = 0.5038461538461538
$ws.Range("K17").Value = 0.08076923076923077
$ws.Range("M17").Value = 0.009615384615384616
$ws.Range("N17").Value = 0.003846153846153846
$ws.Range("O17").Value = 0.09038461538461538
$ws.Range("S17").Value = 0.08269230769230769
$ws.Range("F18").Value = 0.02777777777777778
$ws.Range("H18").Value = 0.1388888888888889
$ws.Range("I18").Value = 0.04166666666666666
$ws.Range("J18").Value = 0.5324074074074074
$ws.Range("K18").Value = 0.05092592592592592
$ws.Range("M18").Value = 0.009259259259259259
$ws.Range("O18").Value = 0.07870370370370371
$ws.Range("S18").Value = 0.1203703703703704
$ws.Range("F19").Value = 0.01970865467009426
$ws.Range("H19").Value = 0.1928020565552699
$ws.Range("I19").Value = 0.03941730934018852
$ws.Range("J19").Value = 0.3933161953727506
$ws.Range("K19").Value = 0.1131105398457583
$ws.Range("M19").Value = 0.02056555269922879
$ws.Range("O19").Value = 0.0908311910882605
$ws.Range("S19").Value = 0.130248500428449
